# Scheduled market-data refresh: update per-Leve price/profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) on each crafting-job sheet
# with freshly scraped values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 22000
$ws.Range("J21").Value = 15000
$ws.Range("L21").Value = 15000
$ws.Range("N21").Value = -15936
$ws.Range("H23").Value = 22000
$ws.Range("J23").Value = 15000
$ws.Range("L23").Value = 15000
$ws.Range("N23").Value = -15468
$ws.Range("H34").Value = 924247.6
$ws.Range("I34").Value = 1118524.9
$ws.Range("K34").Value = 1118524.9
$ws.Range("M34").Value = -1118321.9
$ws.Range("H36").Value = 924247.6
$ws.Range("I36").Value = 1118524.9
$ws.Range("K36").Value = 1118524.9
$ws.Range("M36").Value = -1117809.9
$ws.Range("H55").Value = 149.55556
$ws.Range("I55").Value = 90.666664
$ws.Range("K55").Value = 90.666664
$ws.Range("M55").Value = 123.333336
$ws.Range("H111").Value = 1033
$ws.Range("I111").Value = 400
$ws.Range("J111").Value = 1666
$ws.Range("K111").Value = 1200
$ws.Range("L111").Value = 4998
$ws.Range("M111").Value = 1867
$ws.Range("N111").Value = -11132
$ws.Range("H125").Value = 37370684
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 37370684
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 336336156
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -336341076
$ws.Range("H132").Value = 41555.52
$ws.Range("I132").Value = 43245.375
$ws.Range("J132").Value = 999
$ws.Range("K132").Value = 129736.125
$ws.Range("L132").Value = 2997
$ws.Range("M132").Value = -127206.125
$ws.Range("N132").Value = -8057
$ws.Range("H135").Value = 2638.1304
$ws.Range("I135").Value = 2079.8572
$ws.Range("J135").Value = 8500
$ws.Range("K135").Value = 18718.7148
$ws.Range("L135").Value = 76500
$ws.Range("M135").Value = -16183.7148
$ws.Range("N135").Value = -81570
$ws.Range("H138").Value = 11515484
$ws.Range("I138").Value = 3290752.5
$ws.Range("J138").Value = 15627850
$ws.Range("K138").Value = 9872257.5
$ws.Range("L138").Value = 46883550
$ws.Range("M138").Value = -9867117.5
$ws.Range("N138").Value = -46893830

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 746.5
$ws.Range("I45").Value = 728.6667
$ws.Range("J45").Value = 1014
$ws.Range("K45").Value = 728.6667
$ws.Range("L45").Value = 1014
$ws.Range("M45").Value = -351.6667
$ws.Range("N45").Value = -1768
$ws.Range("H61").Value = 4083.7297
$ws.Range("I61").Value = 3644.5
$ws.Range("K61").Value = 3644.5
$ws.Range("M61").Value = -3432.5
$ws.Range("H96").Value = 26344
$ws.Range("J96").Value = 26344
$ws.Range("L96").Value = 26344
$ws.Range("N96").Value = -31836
$ws.Range("H97").Value = 14977.143
$ws.Range("I97").Value = 17356.666
$ws.Range("K97").Value = 17356.666
$ws.Range("M97").Value = -16860.666
$ws.Range("H110").Value = 1004.14813
$ws.Range("I110").Value = 777.2222
$ws.Range("K110").Value = 777.2222
$ws.Range("M110").Value = 1267.7778
$ws.Range("H132").Value = 3032.0833
$ws.Range("I132").Value = 2585.75
$ws.Range("J132").Value = 4594.25
$ws.Range("K132").Value = 7757.25
$ws.Range("L132").Value = 13782.75
$ws.Range("M132").Value = -5227.25
$ws.Range("N132").Value = -18842.75
$ws.Range("H136").Value = 4083.7297
$ws.Range("I136").Value = 3644.5
$ws.Range("K136").Value = 10933.5
$ws.Range("M136").Value = -8383.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 280748.62
$ws.Range("I105").Value = 2728.077
$ws.Range("J105").Value = 1003602.1
$ws.Range("K105").Value = 2728.077
$ws.Range("L105").Value = 1003602.1
$ws.Range("M105").Value = -981.0770000000002
$ws.Range("N105").Value = -1007096.1
$ws.Range("H134").Value = 2739.6924
$ws.Range("I134").Value = 1528.8334
$ws.Range("J134").Value = 5464.125
$ws.Range("K134").Value = 4586.5002
$ws.Range("L134").Value = 16392.375
$ws.Range("M134").Value = -2051.5002
$ws.Range("N134").Value = -21462.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 985.8889
$ws.Range("I16").Value = 911.2308
$ws.Range("J16").Value = 1180
$ws.Range("K16").Value = 911.2308
$ws.Range("L16").Value = 1180
$ws.Range("M16").Value = -624.2308
$ws.Range("N16").Value = -1754
$ws.Range("H22").Value = 359.66666
$ws.Range("I22").Value = 279.5
$ws.Range("J22").Value = 520
$ws.Range("K22").Value = 279.5
$ws.Range("L22").Value = 520
$ws.Range("M22").Value = 70.5
$ws.Range("N22").Value = -1220
$ws.Range("H58").Value = 2538.8708
$ws.Range("I58").Value = 1137.75
$ws.Range("K58").Value = 1137.75
$ws.Range("M58").Value = -934.75
$ws.Range("H105").Value = 870.5714
$ws.Range("I105").Value = 775.55554
$ws.Range("J105").Value = 1440.6666
$ws.Range("K105").Value = 775.55554
$ws.Range("L105").Value = 1440.6666
$ws.Range("M105").Value = 971.44446
$ws.Range("N105").Value = -4934.6666
$ws.Range("H113").Value = 985.8889
$ws.Range("I113").Value = 911.2308
$ws.Range("J113").Value = 1180
$ws.Range("K113").Value = 911.2308
$ws.Range("L113").Value = 1180
$ws.Range("M113").Value = 1258.7692
$ws.Range("N113").Value = -5520
$ws.Range("H132").Value = 4418.5625
$ws.Range("I132").Value = 4172.636
$ws.Range("J132").Value = 4959.6
$ws.Range("K132").Value = 12517.908
$ws.Range("L132").Value = 14878.8
$ws.Range("M132").Value = -9987.908000000001
$ws.Range("N132").Value = -19938.8
$ws.Range("H134").Value = 2673.625
$ws.Range("I134").Value = 1226.0952
$ws.Range("J134").Value = 5437.091
$ws.Range("K134").Value = 3678.2856
$ws.Range("L134").Value = 16311.273
$ws.Range("M134").Value = -1143.2856
$ws.Range("N134").Value = -21381.273
$ws.Range("H135").Value = 40922
$ws.Range("J135").Value = 40922
$ws.Range("L135").Value = 40922
$ws.Range("N135").Value = -51062
$ws.Range("H136").Value = 2538.8708
$ws.Range("I136").Value = 1137.75
$ws.Range("K136").Value = 3413.25
$ws.Range("M136").Value = -863.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 1186.1428
$ws.Range("I80").Value = 750
$ws.Range("J80").Value = 1360.6
$ws.Range("K80").Value = 2250
$ws.Range("L80").Value = 4081.8
$ws.Range("M80").Value = -1314
$ws.Range("N80").Value = -5953.799999999999
$ws.Range("H83").Value = 1186.1428
$ws.Range("I83").Value = 750
$ws.Range("J83").Value = 1360.6
$ws.Range("K83").Value = 6750
$ws.Range("L83").Value = 12245.4
$ws.Range("M83").Value = -2070
$ws.Range("N83").Value = -21605.4
$ws.Range("H100").Value = 4975
$ws.Range("J100").Value = 4975
$ws.Range("L100").Value = 14925
$ws.Range("N100").Value = -16547

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2761.6
$ws.Range("I132").Value = 2730.8462
$ws.Range("J132").Value = 2794.9167
$ws.Range("K132").Value = 8192.5386
$ws.Range("L132").Value = 8384.750100000001
$ws.Range("M132").Value = -5662.5386
$ws.Range("N132").Value = -13444.7501

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4864.385
$ws.Range("I132").Value = 2670.1
$ws.Range("K132").Value = 8010.299999999999
$ws.Range("M132").Value = -5480.299999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 936
$ws.Range("I96").Value = 800
$ws.Range("J96").Value = 1004
$ws.Range("K96").Value = 800
$ws.Range("L96").Value = 1004
$ws.Range("M96").Value = 573
$ws.Range("N96").Value = -3750
$ws.Range("H132").Value = 1996.6072
$ws.Range("I132").Value = 2038.8383
$ws.Range("K132").Value = 6116.5149
$ws.Range("M132").Value = -3586.5149
$ws.Range("H136").Value = 1180.38
$ws.Range("I136").Value = 662.8484999999999
$ws.Range("J136").Value = 2185
$ws.Range("K136").Value = 1988.5455
$ws.Range("L136").Value = 6555
$ws.Range("M136").Value = 561.4545000000003
$ws.Range("N136").Value = -11655
